# "commit on 15th may with config reader"
#
# Appends four new RegisterPage test rows (8-11) plus a final "account
# created" row (12) to dsAlgoInput.xlsx, wires up mailto hyperlinks for the
# new UserName/Password cells (matching the sheet's existing convention),
# and widens column E so the longer messages stay readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: tag8 -------------------------------------------------------
$ws.Range("A8").Value = "tag8"
$ws.Range("B8").Value = "Numpy@sdet84_1"
$ws.Range("C8").Value = "pass"
$ws.Range("D8").Value = "pass"
$ws.Range("E8").Value = "Password should contain at least 8 characters"

# --- Row 9: tag9 (password that is entirely numeric) -------------------
$ws.Range("E9").Value = "Password can’t be entirely numeric."
$ws.Range("A9").Value = "tag9"
$ws.Range("B9").Value = "Numpy@sdet84_1"
$ws.Range("C9").Value = 123456789
$ws.Range("D9").Value = 123456789

# --- Tags for rows 10-12 -------------------------------------------------
$ws.Range("A10").Value = "tag10"
$ws.Range("A11").Value = "tag11"
$ws.Range("A12").Value = "tag12"

# --- Row 10 remainder: tag10 (password too similar to personal info) ---
$ws.Range("B10").Value = "Numpy@sdet84_1"
$ws.Range("C10").Value = "testsdet84"
$ws.Range("D10").Value = "testsdet84"
$ws.Range("E10").Value = "password can’t be too similar to your other personal information."

# --- Row 11 remainder: tag11 (commonly used password) ------------------
$ws.Range("B11").Value = "Numpy@sdet84_1"
$ws.Range("C11").Value = "welcome"
$ws.Range("D11").Value = "welcome"
$ws.Range("E11").Value = "Password can’t be commonly used password"

# --- Row 12 remainder: tag12 (successful account creation) -------------
$ws.Range("B12").Value = "numpyqueen"
$ws.Range("C12").Value = "queen@1305"
$ws.Range("D12").Value = "queen@1305"
$ws.Range("E12").Value = "New Account Created. You are logged in as numpyqueen"

# --- mailto hyperlinks for the new UserName / Password cells -----------
# (Hyperlinks.Add stamps the built-in "Hyperlink" style on the target cell,
# so strip it back off for the cells that should stay plain.)
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:Numpy@sdet84_1")
$ws.Range("B8").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:Numpy@sdet84_1")
$ws.Range("B9").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:Numpy@sdet84_1")
$ws.Range("B10").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Numpy@sdet84_1")
$ws.Range("B11").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:queen@1305")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:queen@1305")

# B7's hyperlink pre-dates this commit; it now picks up the same
# Hyperlink visual style as the other linked cells.
$ws.Range("B7").Style = "Hyperlink"

# --- Column E widens to fit the new, longer messages --------------------
$ws.Columns("E").ColumnWidth = 84

# --- Selection follows the new last cell, like Excel leaves behind ------
$ws.Range("E12").Select()
